$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting rows 65:183 down to 66:184
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record's data
$ws.Cells.Item(65, 1).Value = 10
$ws.Cells.Item(65, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(65, 3).Value = "La Araucanía"
$ws.Cells.Item(65, 4).Value = 45002
$ws.Cells.Item(65, 5).Value = 9
$ws.Cells.Item(65, 6).Value = 100114002
$ws.Cells.Item(65, 7).Value = "Camote"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 50
$ws.Cells.Item(65, 11).Value = 26000
$ws.Cells.Item(65, 12).Value = 26000
$ws.Cells.Item(65, 13).Value = 26000
$ws.Cells.Item(65, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(65, 15).Value = "Perú"
$ws.Cells.Item(65, 16).Value = 1300
$ws.Cells.Item(65, 17).Value = 20
$ws.Cells.Item(65, 18).Value = "Hortaliza"

# Apply the same style (date format) as used in column D for the rest of the table
$ws.Cells.Item(65, 4).NumberFormat = $ws.Cells.Item(66, 4).NumberFormat
